$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "algemeen" sheet: add swing_output / swing_output_bestandsnaam columns
# ---------------------------------------------------------------------------
$algemeen = $wb.Worksheets.Item("algemeen")
$algemeen.Range("S1").Value = "swing_output"
$algemeen.Range("T1").Value = "swing_output_bestandsnaam"
$algemeen.Range("S2").Value = $true
$algemeen.Range("T2").Value = "kubusdata"

# ---------------------------------------------------------------------------
# 2. Add the new "swing" worksheet, positioned right before "logos"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "swing"
$logos = $wb.Worksheets.Item("logos")
$newSheet.Move($logos)

$swing = $wb.Worksheets.Item("swing")

$swing.Range("A1").Value = "inidcator_val"
$swing.Range("B1").Value = "indicator_code"
$swing.Range("C1").Value = "indicator_label"
$swing.Range("D1").Value = "measure"
$swing.Range("E1").Value = "unit"
$swing.Range("F1").Value = "decimals"
$swing.Range("G1").Value = "geo_subset"
$swing.Range("H1").Value = "geo_code_source"
$swing.Range("I1").Value = "include_crossing"

$swing.Range("A2").Value = "SLAAPSLECHT"
$swing.Range("B2").Value = 1
$swing.Range("D2").Value = "perc"
$swing.Range("E2").Value = "%"
$swing.Range("I2").Value = "TRUE"

$swing.Columns.Item(1).ColumnWidth = 12.14

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping to mirror the authored workbook
# ---------------------------------------------------------------------------
$onderdelen = $wb.Worksheets.Item("onderdelen")
$onderdelen.Range("B1").Select()

$opmaak = $wb.Worksheets.Item("opmaak")
$opmaak.Range("B29").Select()

$indelingRijen = $wb.Worksheets.Item("indeling_rijen")
$indelingRijen.Range("E7").Select()

$algemeen2 = $wb.Worksheets.Item("algemeen")
$algemeen2.Range("P23").Select()

$wb.Worksheets.Item("swing").Activate()
$wb.Worksheets.Item("swing").Range("J2").Select()
